# 艺术表演场馆艺术演出场次.xlsx — update data series:
#   - drop the oldest two years (2008年, 2009年) -> rows 2 and 3
#   - every remaining year row shifts up by two rows
#   - append a new latest year (2021年) as the new last data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2008年 / 2009年 rows (rows 2 and 3); rows below shift up by 2,
# so what used to be 2010年 (row 4) becomes the new row 2, etc.
$ws.Range("2:3").EntireRow.Delete()

# After the deletion, the last data row is 2020年 on row 12 and the sheet's
# used range is A1:S12. Add the new 2021年 row right after it as row 13.
$lastRow = 12
$newRow = $lastRow + 1

# Match the formatting (bold/centered/bordered year-label style) used by the
# other year cells in column A by copying it from the row above.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A" + $newRow).Value = "2021年"
$ws.Range("B" + $newRow).Value = 0.5
$ws.Range("C" + $newRow).Value = 9.199999999999999
$ws.Range("D" + $newRow).Value = 59.5
# Column E ("其他部门艺术表演场馆艺术演出场次") stays blank for this row,
# matching the rest of the recent-years series.
$ws.Range("F" + $newRow).Value = 105.1
$ws.Range("G" + $newRow).Value = 592.9
$ws.Range("H" + $newRow).Value = 60.1
$ws.Range("I" + $newRow).Value = 36.3
$ws.Range("J" + $newRow).Value = 404.6
$ws.Range("K" + $newRow).Value = 70.7
$ws.Range("L" + $newRow).Value = 571.9
# Column M ("文化部门艺术表演场馆艺术演出场次") also stays blank here.
$ws.Range("N" + $newRow).Value = 4.7
$ws.Range("O" + $newRow).Value = 12.9
$ws.Range("P" + $newRow).Value = 43.4
$ws.Range("Q" + $newRow).Value = 642.6
$ws.Range("R" + $newRow).Value = 1.6
$ws.Range("S" + $newRow).Value = 16.1
